$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), columns A..T
$headers = @("ID", "Family", "Genus", "Species", "Sex", "Pinned", "Body", "Forewing.left", "Forewing.right", "Hindwing.left", "Hindwing.right", "Site", "Latitude", "Longitude", "Climate", "Date", "Collector", "Binomial", "DNA", "Spectra")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Row 2
$row2 = @(4, "Hesperiidae", "Notocrypta", "waigensis", "Female", "n", "n", "n", "n", "n", "y", "BG", -16.90136852, 145.7505777, "Tropical", "16/07/2022", "MFE", "Notocrypta waigensis", "n", "n")
for ($c = 1; $c -le $row2.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2[$c - 1]
}

# Row 3
$row3 = @(16, "Hesperiidae", "Telicota", "mesoptis", "Male", "n", "n", "n", "n", "n", "n", "BG", -16.90136852, 145.7505777, "Tropical", "16/07/2022", "MFE", "Telicota mesoptis", "n", "n")
for ($c = 1; $c -le $row3.Length; $c++) {
    $ws.Cells.Item(3, $c).Value = $row3[$c - 1]
}

# Row 4
$row4 = @(19, "Hesperiidae", "Telicota", "mesoptis", "Male", "n", "n", "n", "n", "n", "n", "BG", -16.90136852, 145.7505777, "Tropical", "16/07/2022", "MFE", "Telicota mesoptis", "n", "n")
for ($c = 1; $c -le $row4.Length; $c++) {
    $ws.Cells.Item(4, $c).Value = $row4[$c - 1]
}

# Row 5
$row5 = @(186, "Hesperiidae", "Suniana", "sunias", "Female", "n", "y", "n", "n", "n", "n", "JCU", -16.81425603, 145.6854526, "Tropical", "16/07/2022", "SS", "Suniana sunias", "y", "y")
for ($c = 1; $c -le $row5.Length; $c++) {
    $ws.Cells.Item(5, $c).Value = $row5[$c - 1]
}

# Row 6
$row6 = @(551, "Nymphalidae", "Euploea", "darchia", "Female", "y", "n", "n", "n", "n", "n", "MR", -16.45193583, 145.3714342, "Tropical", "19/07/2022", "SS", "Euploea darchia", "y", "y")
for ($c = 1; $c -le $row6.Length; $c++) {
    $ws.Cells.Item(6, $c).Value = $row6[$c - 1]
}

# Row 7
$row7 = @(1361, "Papilionidae", "Papilio", "aegeus", "Male", "n", "n", "n", "n", "n", "n", "CC", -27.344776, 153.0390842, "Subtropical", "24/09/2022", "CI", "Papilio aegeus", "n", "n")
for ($c = 1; $c -le $row7.Length; $c++) {
    $ws.Cells.Item(7, $c).Value = $row7[$c - 1]
}
